$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37707
$ws.Range("D2").Value = 54533235
$ws.Range("C3").Value = 90917
$ws.Range("D3").Value = 133277715
$ws.Range("C4").Value = 31163
$ws.Range("D4").Value = 46150849
$ws.Range("C5").Value = 8689
$ws.Range("D5").Value = 12914063
$ws.Range("C6").Value = 1991
$ws.Range("D6").Value = 2959006
$ws.Range("C12").Value = 41292
$ws.Range("D12").Value = 56027689
$ws.Range("C13").Value = 9648
$ws.Range("D13").Value = 13953958
$ws.Range("C14").Value = 25933
$ws.Range("D14").Value = 38035540
$ws.Range("C15").Value = 8308
$ws.Range("D15").Value = 12329824
$ws.Range("C16").Value = 2151
$ws.Range("D16").Value = 3198665
$ws.Range("C20").Value = 10213
$ws.Range("D20").Value = 13525214
$ws.Range("C21").Value = 13372
$ws.Range("D21").Value = 19308292
$ws.Range("C22").Value = 31645
$ws.Range("D22").Value = 46439362
$ws.Range("C23").Value = 10214
$ws.Range("D23").Value = 15183678
$ws.Range("C27").Value = 11680
$ws.Range("D27").Value = 15602344
$ws.Range("C28").Value = 7636
$ws.Range("D28").Value = 11060617
$ws.Range("C29").Value = 22471
$ws.Range("D29").Value = 32984549
$ws.Range("C30").Value = 7811
$ws.Range("D30").Value = 11624133
$ws.Range("C34").Value = 8302
$ws.Range("D34").Value = 10966515
$ws.Range("C35").Value = 3243
$ws.Range("D35").Value = 4681194
$ws.Range("C36").Value = 7820
$ws.Range("D36").Value = 11419934
$ws.Range("C37").Value = 3176
$ws.Range("D37").Value = 4706961
$ws.Range("C38").Value = 829
$ws.Range("D38").Value = 1234723
$ws.Range("C39").Value = 164
$ws.Range("D39").Value = 243686
$ws.Range("C42").Value = 17222
$ws.Range("D42").Value = 24902458
$ws.Range("C43").Value = 51066
$ws.Range("D43").Value = 74862786
$ws.Range("C44").Value = 19005
$ws.Range("D44").Value = 28229943
$ws.Range("C45").Value = 5604
$ws.Range("D45").Value = 8345177
$ws.Range("C46").Value = 1202
$ws.Range("D46").Value = 1793545
$ws.Range("C50").Value = 16686
$ws.Range("D50").Value = 22216188
$ws.Range("C51").Value = 2015
$ws.Range("D51").Value = 2922108
$ws.Range("C52").Value = 6885
$ws.Range("D52").Value = 10121074
$ws.Range("C53").Value = 2345
$ws.Range("D53").Value = 3502418
$ws.Range("C57").Value = 6953
$ws.Range("D57").Value = 9561256
$ws.Range("C58").Value = 943
$ws.Range("D58").Value = 1384079
$ws.Range("C59").Value = 2369
$ws.Range("D59").Value = 3512337
$ws.Range("C60").Value = 940
$ws.Range("D60").Value = 1399501
$ws.Range("C64").Value = 1389
$ws.Range("D64").Value = 1954706
$ws.Range("C65").Value = 15341
$ws.Range("D65").Value = 22161331
$ws.Range("C66").Value = 44650
$ws.Range("D66").Value = 65339707
$ws.Range("C67").Value = 15690
$ws.Range("D67").Value = 23317186
$ws.Range("C68").Value = 4566
$ws.Range("D68").Value = 6801292
$ws.Range("C69").Value = 923
$ws.Range("D69").Value = 1372668
$ws.Range("C73").Value = 15070
$ws.Range("D73").Value = 19871392
$ws.Range("C74").Value = 51313
$ws.Range("D74").Value = 74672749
$ws.Range("C75").Value = 145861
$ws.Range("D75").Value = 214889029
$ws.Range("C76").Value = 63565
$ws.Range("D76").Value = 94720615
$ws.Range("C77").Value = 20320
$ws.Range("D77").Value = 30360331
$ws.Range("C78").Value = 4810
$ws.Range("D78").Value = 7184043
$ws.Range("C81").Value = 15
$ws.Range("D81").Value = 21625
$ws.Range("C85").Value = 50746
$ws.Range("D85").Value = 69034703
$ws.Range("C86").Value = 4593
$ws.Range("D86").Value = 6654936
$ws.Range("C87").Value = 11550
$ws.Range("D87").Value = 16968542
$ws.Range("C93").Value = 5404
$ws.Range("D93").Value = 7264956
$ws.Range("C94").Value = 1593
$ws.Range("D94").Value = 2294432
$ws.Range("C95").Value = 5155
$ws.Range("D95").Value = 7591243
$ws.Range("C97").Value = 689
$ws.Range("D97").Value = 1032460
$ws.Range("C101").Value = 3550
$ws.Range("D101").Value = 4697961
$ws.Range("C103").Value = 352
$ws.Range("D103").Value = 525530
$ws.Range("C107").Value = 10744
$ws.Range("D107").Value = 15585962
$ws.Range("C108").Value = 29179
$ws.Range("D108").Value = 42871096
$ws.Range("C109").Value = 9772
$ws.Range("D109").Value = 14531650
$ws.Range("C114").Value = 9784
$ws.Range("D114").Value = 12926162
$ws.Range("C115").Value = 30403
$ws.Range("D115").Value = 43843724
$ws.Range("C116").Value = 66093
$ws.Range("D116").Value = 96725176
$ws.Range("C117").Value = 21352
$ws.Range("D117").Value = 31732440
$ws.Range("C118").Value = 6060
$ws.Range("D118").Value = 9028021
$ws.Range("C124").Value = 25816
$ws.Range("D124").Value = 34484907
$ws.Range("C125").Value = 35944
$ws.Range("D125").Value = 51875027
$ws.Range("C126").Value = 76707
$ws.Range("D126").Value = 112168093
$ws.Range("C127").Value = 23825
$ws.Range("D127").Value = 35358909
$ws.Range("C128").Value = 6387
$ws.Range("D128").Value = 9491238
$ws.Range("C129").Value = 1233
$ws.Range("D129").Value = 1833911
$ws.Range("C133").Value = 31770
$ws.Range("D133").Value = 42186614
$ws.Range("C134").Value = 13199
$ws.Range("D134").Value = 19104630
$ws.Range("C135").Value = 32299
$ws.Range("D135").Value = 47439783
$ws.Range("C136").Value = 11463
$ws.Range("D136").Value = 17032542
$ws.Range("C138").Value = 499
$ws.Range("D138").Value = 742490
$ws.Range("C141").Value = 10801
$ws.Range("D141").Value = 14404056
$ws.Range("C142").Value = 34988
$ws.Range("D142").Value = 50524142
$ws.Range("C143").Value = 81137
$ws.Range("D143").Value = 118875615
$ws.Range("C144").Value = 24314
$ws.Range("D144").Value = 36125155
$ws.Range("C149").Value = 29151
$ws.Range("D149").Value = 39327767
